$wb = $excel.ActiveWorkbook

# ALC row 11: Gotta Bounce
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 77340.7
$ws.Range("I11").Value = 77340.7
$ws.Range("K11").Value = 77340.7
$ws.Range("M11").Value = -77200.7

# ALC row 42: Eye of the Beholder
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 62
$ws.Range("I42").Value = 64.44444
$ws.Range("K42").Value = 193.33332
$ws.Range("M42").Value = 36.66667999999999

# ALC row 76: Warding Off Temptation
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4252.0835
$ws.Range("I76").Value = 3376.125
$ws.Range("K76").Value = 3376.125
$ws.Range("M76").Value = -3061.125

# ALC row 79: The Garden of Arcane Delights (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4252.0835
$ws.Range("I79").Value = 3376.125
$ws.Range("K79").Value = 3376.125
$ws.Range("M79").Value = -2284.125

# ALC row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 34291.168
$ws.Range("I132").Value = 50745.75
$ws.Range("K132").Value = 152237.25
$ws.Range("M132").Value = -149707.25

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 177058.62
$ws.Range("I32").Value = 178410.53
$ws.Range("K32").Value = 178410.53
$ws.Range("M32").Value = -178123.53

# ARM row 45: Hollow Hallmarks
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4062.5
$ws.Range("I45").Value = 3750
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 3750
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -3373
$ws.Range("N45").Value = -5754

# ARM row 74: As the Bolt Flies
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5652.8115
$ws.Range("I74").Value = 2692.4092
$ws.Range("K74").Value = 2692.4092
$ws.Range("M74").Value = -1818.4092

# ARM row 77: Heavy Metal Banned (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5652.8115
$ws.Range("I77").Value = 2692.4092
$ws.Range("K77").Value = 13462.046
$ws.Range("M77").Value = -9094.046

# ARM row 80: A Squire to Inspire
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 37500
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# ARM row 83: All's Fair in Highborn Assassination (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 37500
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# ARM row 110: Scheduled Maintenance
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1038.6522
$ws.Range("I110").Value = 1026.8096
$ws.Range("K110").Value = 1026.8096
$ws.Range("M110").Value = 1018.1904

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5266.5894
$ws.Range("I132").Value = 3367.697
$ws.Range("K132").Value = 10103.091
$ws.Range("M132").Value = -7573.091

# BSM row 20: Smelt and Dealt
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 46753.176
$ws.Range("J20").Value = 2430.7
$ws.Range("L20").Value = 2430.7
$ws.Range("N20").Value = -2924.7

# BSM row 99: Meddle in Metal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 9966.333000000001
$ws.Range("J99").Value = 1699.3334
$ws.Range("L99").Value = 1699.3334
$ws.Range("N99").Value = -4695.3334

# CRP row 14: Citizens' Canes
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 5236.6
$ws.Range("I14").Value = 4750
$ws.Range("K14").Value = 4750
$ws.Range("M14").Value = -4580

# CRP row 25: Bowing to Necessity
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 9511
$ws.Range("I25").Value = 9511
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 9511
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -9337
$ws.Range("N25").ClearContents()

# CRP row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2148
$ws.Range("I31").Value = 2826.2727
$ws.Range("K31").Value = 2826.2727
$ws.Range("M31").Value = -2531.2727

# CRP row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2148
$ws.Range("I34").Value = 2826.2727
$ws.Range("K34").Value = 2826.2727
$ws.Range("M34").Value = -2624.2727

# CRP row 58: You Do the Heavy Lifting
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5923.9443
$ws.Range("I58").Value = 4390.846
$ws.Range("J58").Value = 9910
$ws.Range("K58").Value = 4390.846
$ws.Range("L58").Value = 9910
$ws.Range("M58").Value = -4187.846
$ws.Range("N58").Value = -10316

# CRP row 94: Beech, Please
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6072.3
$ws.Range("J94").Value = 1408.6923
$ws.Range("L94").Value = 1408.6923
$ws.Range("N94").Value = -2310.6923

# CRP row 105: Zelkova, My Love
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 9011.429
$ws.Range("I105").Value = 9550.77
$ws.Range("K105").Value = 9550.77
$ws.Range("M105").Value = -7803.77

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3316.75
$ws.Range("I132").Value = 3359.7368
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 10079.2104
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -7549.2104
$ws.Range("N132").Value = -12560

# CRP row 134: Wood You Be Quiet
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2291.524
$ws.Range("I134").Value = 2291.524
$ws.Range("K134").Value = 6874.572
$ws.Range("M134").Value = -4339.572

# CRP row 136: Turali Quality
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5923.9443
$ws.Range("I136").Value = 4390.846
$ws.Range("J136").Value = 9910
$ws.Range("K136").Value = 13172.538
$ws.Range("L136").Value = 29730
$ws.Range("M136").Value = -10622.538
$ws.Range("N136").Value = -34830

# CUL row 12: Butter Me Up
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 371.3846
$ws.Range("I12").Value = 950
$ws.Range("J12").Value = 323.16666
$ws.Range("K12").Value = 2850
$ws.Range("L12").Value = 969.4999799999999
$ws.Range("M12").Value = -2677
$ws.Range("N12").Value = -1315.49998

# CUL row 37: I Love Lamprey
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 57800
$ws.Range("J37").Value = 57800
$ws.Range("L37").Value = 173400
$ws.Range("N37").Value = -173624

# GSM row 70: Sky Is the Limit
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18331.268
$ws.Range("I70").Value = 28125
$ws.Range("J70").Value = 7138.4287
$ws.Range("K70").Value = 28125
$ws.Range("L70").Value = 7138.4287
$ws.Range("M70").Value = -27855
$ws.Range("N70").Value = -7678.4287

# GSM row 73: Hulls of Broken Dreams (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 18331.268
$ws.Range("I73").Value = 28125
$ws.Range("J73").Value = 7138.4287
$ws.Range("K73").Value = 28125
$ws.Range("L73").Value = 7138.4287
$ws.Range("M73").Value = -27189
$ws.Range("N73").Value = -9010.4287

# GSM row 132: On Board for Lar
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 13983.909
$ws.Range("J132").Value = 3659.8
$ws.Range("L132").Value = 10979.4
$ws.Range("N132").Value = -16039.4

# LTW row 46: Supply Side Logic
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4439.6924
$ws.Range("J46").Value = 4581.7915
$ws.Range("L46").Value = 4581.7915
$ws.Range("N46").Value = -4957.7915

# LTW row 61: Spelling Me Softly
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7082.6587
$ws.Range("I61").Value = 6912.4
$ws.Range("K61").Value = 6912.4
$ws.Range("M61").Value = -6710.4

# LTW row 100: Tiger in the Sack
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3844.7778
$ws.Range("I100").Value = 3941.4
$ws.Range("J100").Value = 3724
$ws.Range("K100").Value = 3941.4
$ws.Range("L100").Value = 3724
$ws.Range("M100").Value = -3400.4
$ws.Range("N100").Value = -4806

# LTW row 113: Peace in Rest
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 7082.6587
$ws.Range("I113").Value = 6912.4
$ws.Range("K113").Value = 6912.4
$ws.Range("M113").Value = -4742.4

# LTW row 122: Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4691.4375
$ws.Range("I122").Value = 3238.5454
$ws.Range("K122").Value = 9715.636200000001
$ws.Range("M122").Value = -7265.636200000001

# LTW row 136: Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2509.6785
$ws.Range("I136").Value = 2142.2173
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 6426.651899999999
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -3876.651899999999
$ws.Range("N136").Value = -17700

# WVR row 81: Where the Dragonflies, the Net Catches
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3481.5
$ws.Range("I81").Value = 2520.8333
$ws.Range("K81").Value = 5041.6666
$ws.Range("M81").Value = -3980.6666

# WVR row 84: To Kill a Dragon on Nameday (L)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3481.5
$ws.Range("I84").Value = 2520.8333
$ws.Range("K84").Value = 25208.333
$ws.Range("M84").Value = -19904.333

# WVR row 132: Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4261.091
$ws.Range("I132").Value = 3971.6
$ws.Range("K132").Value = 11914.8
$ws.Range("M132").Value = -9384.799999999999
